# Auto-generated script applying scheduled market-price refresh to Jenova_Profits sheets
$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 799.4
$ws.Range("J18").Value = 999
$ws.Range("L18").Value = 999
$ws.Range("N18").Value = -1567
$ws.Range("H34").Value = 22261
$ws.Range("I34").Value = 22261
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 22261
$ws.Range("L34").Value = 0
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -22058
$ws.Range("H36").Value = 22261
$ws.Range("I36").Value = 22261
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 22261
$ws.Range("L36").Value = 0
$ws.Range("M36").ClearContents()
$ws.Range("N36").Value = -21546
$ws.Range("H44").Value = 376870.34
$ws.Range("J44").Value = 9500
$ws.Range("L44").Value = 9500
$ws.Range("N44").Value = -10424
$ws.Range("H74").Value = 10912.714
$ws.Range("I74").Value = 8981.5
$ws.Range("K74").Value = 8981.5
$ws.Range("M74").Value = -8045.5
$ws.Range("H77").Value = 10912.714
$ws.Range("I77").Value = 8981.5
$ws.Range("K77").Value = 44907.5
$ws.Range("M77").Value = -40227.5
$ws.Range("H99").Value = 708.1818
$ws.Range("J99").Value = 2023
$ws.Range("L99").Value = 6069
$ws.Range("N99").Value = -9065
$ws.Range("H111").Value = 132759.5
$ws.Range("I111").Value = 206811
$ws.Range("J111").Value = 9340.333000000001
$ws.Range("K111").Value = 620433
$ws.Range("L111").Value = 28020.999
$ws.Range("M111").Value = -617366
$ws.Range("N111").Value = -34154.999
$ws.Range("H132").Value = 1399.275
$ws.Range("I132").Value = 1515.8572
$ws.Range("K132").Value = 4547.571599999999
$ws.Range("M132").Value = -2017.571599999999
$ws.Range("H136").Value = 105000
$ws.Range("J136").Value = 105000
$ws.Range("L136").Value = 105000
$ws.Range("N136").Value = -115200
$ws.Range("H137").Value = 2807.4092
$ws.Range("J137").Value = 3330.1667
$ws.Range("L137").Value = 9990.500100000001
$ws.Range("N137").Value = -15090.5001
$ws.Range("H138").Value = 4431.909
$ws.Range("J138").Value = 5752.931
$ws.Range("L138").Value = 17258.793
$ws.Range("N138").Value = -27538.793

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4663.9077
$ws.Range("I32").Value = 4663.9077
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 4663.9077
$ws.Range("L32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -4376.9077
$ws.Range("H132").Value = 2518.4187
$ws.Range("I132").Value = 2605.639
$ws.Range("J132").Value = 2069.8572
$ws.Range("K132").Value = 7816.917
$ws.Range("L132").Value = 6209.571599999999
$ws.Range("M132").Value = -5286.917
$ws.Range("N132").Value = -11269.5716

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 16890
$ws.Range("I26").Value = 16890
$ws.Range("K26").Value = 16890
$ws.Range("M26").Value = -16598
$ws.Range("H108").Value = 80695
$ws.Range("J108").Value = 80695
$ws.Range("L108").Value = 80695
$ws.Range("N108").Value = -88375
$ws.Range("H134").Value = 96835.27
$ws.Range("I134").Value = 6518.8
$ws.Range("K134").Value = 19556.4
$ws.Range("M134").Value = -17021.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 6236.0527
$ws.Range("I16").Value = 2248.7693
$ws.Range("J16").Value = 14875.167
$ws.Range("K16").Value = 2248.7693
$ws.Range("L16").Value = 14875.167
$ws.Range("M16").Value = -1961.7693
$ws.Range("N16").Value = -15449.167
$ws.Range("H62").Value = 3114.625
$ws.Range("I62").Value = 2416.7144
$ws.Range("J62").Value = 8000
$ws.Range("K62").Value = 2416.7144
$ws.Range("L62").Value = 8000
$ws.Range("M62").Value = -1792.7144
$ws.Range("N62").Value = -9248
$ws.Range("H65").Value = 3114.625
$ws.Range("I65").Value = 2416.7144
$ws.Range("J65").Value = 8000
$ws.Range("K65").Value = 12083.572
$ws.Range("L65").Value = 40000
$ws.Range("M65").Value = -8963.572
$ws.Range("N65").Value = -46240
$ws.Range("H99").Value = 388608.06
$ws.Range("I99").Value = 3590.5
$ws.Range("J99").Value = 1672000
$ws.Range("K99").Value = 3590.5
$ws.Range("L99").Value = 1672000
$ws.Range("M99").Value = -2092.5
$ws.Range("N99").Value = -1674996
$ws.Range("H113").Value = 6236.0527
$ws.Range("I113").Value = 2248.7693
$ws.Range("J113").Value = 14875.167
$ws.Range("K113").Value = 2248.7693
$ws.Range("L113").Value = 14875.167
$ws.Range("M113").Value = -78.76929999999993
$ws.Range("N113").Value = -19215.167
$ws.Range("H122").Value = 2984
$ws.Range("I122").Value = 3457.8
$ws.Range("J122").Value = 1799.5
$ws.Range("K122").Value = 10373.4
$ws.Range("L122").Value = 5398.5
$ws.Range("M122").Value = -7923.400000000001
$ws.Range("N122").Value = -10298.5
$ws.Range("H126").Value = 388608.06
$ws.Range("I126").Value = 3590.5
$ws.Range("J126").Value = 1672000
$ws.Range("K126").Value = 10771.5
$ws.Range("L126").Value = 5016000
$ws.Range("M126").Value = -8301.5
$ws.Range("N126").Value = -5020940
$ws.Range("H132").Value = 1165
$ws.Range("I132").Value = 925.125
$ws.Range("J132").Value = 1804.6666
$ws.Range("K132").Value = 2775.375
$ws.Range("L132").Value = 5413.9998
$ws.Range("M132").Value = -245.375
$ws.Range("N132").Value = -10473.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 16000
$ws.Range("I39").Value = 833.3333
$ws.Range("K39").Value = 2499.9999
$ws.Range("M39").Value = -2205.9999
$ws.Range("H60").Value = 508.33334
$ws.Range("I60").Value = 410
$ws.Range("J60").Value = 1000
$ws.Range("K60").Value = 1230
$ws.Range("L60").Value = 3000
$ws.Range("M60").Value = -979
$ws.Range("N60").Value = -3502
$ws.Range("H109").Value = 55900.555
$ws.Range("I109").Value = 365.29413
$ws.Range("K109").Value = 1095.88239
$ws.Range("M109").Value = -55.88238999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 7512500
$ws.Range("J11").Value = 6683333.5
$ws.Range("L11").Value = 6683333.5
$ws.Range("N11").Value = -6683611.5
$ws.Range("H102").Value = 1572.25
$ws.Range("I102").Value = 670.2105
$ws.Range("J102").Value = 5000
$ws.Range("K102").Value = 670.2105
$ws.Range("L102").Value = 5000
$ws.Range("M102").Value = 951.7895
$ws.Range("N102").Value = -8244
$ws.Range("H122").Value = 4092.3572
$ws.Range("I122").Value = 2742
$ws.Range("J122").Value = 5442.7144
$ws.Range("K122").Value = 8226
$ws.Range("L122").Value = 16328.1432
$ws.Range("M122").Value = -5776
$ws.Range("N122").Value = -21228.1432
$ws.Range("H126").Value = 3604.3333
$ws.Range("I126").Value = 3348.6
$ws.Range("J126").Value = 3924
$ws.Range("K126").Value = 10045.8
$ws.Range("L126").Value = 11772
$ws.Range("M126").Value = -7575.799999999999
$ws.Range("N126").Value = -16712
$ws.Range("H140").Value = 68624.875
$ws.Range("J140").Value = 68624.875
$ws.Range("L140").Value = 68624.875
$ws.Range("N140").Value = -78984.875
$ws.Range("H141").Value = 69500
$ws.Range("J141").Value = 69500
$ws.Range("L141").Value = 69500
$ws.Range("N141").Value = -79860

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 6463.25
$ws.Range("I61").Value = 6672.4287
$ws.Range("K61").Value = 6672.4287
$ws.Range("M61").Value = -6470.4287
$ws.Range("H113").Value = 6463.25
$ws.Range("I113").Value = 6672.4287
$ws.Range("K113").Value = 6672.4287
$ws.Range("M113").Value = -4502.4287
$ws.Range("H124").Value = 31334.8
$ws.Range("J124").Value = 31334.8
$ws.Range("L124").Value = 31334.8
$ws.Range("N124").Value = -41154.8
$ws.Range("H136").Value = 1824825.8
$ws.Range("I136").Value = 2506172
$ws.Range("K136").Value = 7518516
$ws.Range("M136").Value = -7515966

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H27").Value = 75310.5
$ws.Range("J27").Value = 75310.5
$ws.Range("L27").Value = 75310.5
$ws.Range("N27").Value = -75448.5
$ws.Range("H56").Value = 5561.6665
$ws.Range("I56").Value = 4000
$ws.Range("J56").Value = 5874
$ws.Range("K56").Value = 4000
$ws.Range("L56").Value = 5874
$ws.Range("M56").Value = -3286
$ws.Range("N56").Value = -7302
$ws.Range("H80").Value = 120000
$ws.Range("J80").Value = 120000
$ws.Range("L80").Value = 120000
$ws.Range("N80").Value = -121996
$ws.Range("H83").Value = 120000
$ws.Range("J83").Value = 120000
$ws.Range("L83").Value = 360000
$ws.Range("N83").Value = -369984
$ws.Range("H132").Value = 112858.78
$ws.Range("I132").Value = 1104.1428
$ws.Range("K132").Value = 3312.4284
$ws.Range("M132").Value = -782.4284000000002
$ws.Range("H135").Value = 99333
$ws.Range("J135").Value = 99333
$ws.Range("L135").Value = 99333
$ws.Range("N135").Value = -109473
$ws.Range("H140").Value = 140000
$ws.Range("J140").Value = 140000
$ws.Range("L140").Value = 140000
$ws.Range("N140").Value = -150360
$ws.Range("H141").Value = 48750
$ws.Range("J141").Value = 48750
$ws.Range("L141").Value = 48750
$ws.Range("N141").Value = -59110
